$p = $ppt.ActivePresentation
$p.Slides.Item(18).Delete()
